$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 218-223 with revised monthly export figures ---
# Row 218
$ws.Cells.Item(218, 2).Value = 7177
$ws.Cells.Item(218, 3).Value = 3647
$ws.Cells.Item(218, 4).Value = 3283
$ws.Cells.Item(218, 5).Value = 1316
$ws.Cells.Item(218, 6).Value = 1855
$ws.Cells.Item(218, 13).Value = 1400
$ws.Cells.Item(218, 14).Value = 1351
$ws.Cells.Item(218, 15).Value = 56
$ws.Cells.Item(218, 18).Value = 201
$ws.Cells.Item(218, 20).Value = 8
$ws.Cells.Item(218, 21).Value = 1012
$ws.Cells.Item(218, 22).Value = 18
$ws.Cells.Item(218, 23).Value = 35
$ws.Cells.Item(218, 28).Value = 2130
$ws.Cells.Item(218, 29).Value = 777
$ws.Cells.Item(218, 32).Value = 393
$ws.Cells.Item(218, 58).Value = 396
$ws.Cells.Item(218, 60).Value = 44
$ws.Cells.Item(218, 63).Value = 75

# Row 219
$ws.Cells.Item(219, 2).Value = 7294
$ws.Cells.Item(219, 3).Value = 4510
$ws.Cells.Item(219, 4).Value = 4254
$ws.Cells.Item(219, 5).Value = 1636
$ws.Cells.Item(219, 6).Value = 2443
$ws.Cells.Item(219, 13).Value = 867
$ws.Cells.Item(219, 14).Value = 825
$ws.Cells.Item(219, 15).Value = 173
$ws.Cells.Item(219, 17).Value = 14
$ws.Cells.Item(219, 18).Value = 198
$ws.Cells.Item(219, 20).Value = 47
$ws.Cells.Item(219, 21).Value = 317
$ws.Cells.Item(219, 22).Value = 8
$ws.Cells.Item(219, 28).Value = 1918
$ws.Cells.Item(219, 29).Value = 745
$ws.Cells.Item(219, 32).Value = 366
$ws.Cells.Item(219, 47).Value = 149
$ws.Cells.Item(219, 48).Value = 48
$ws.Cells.Item(219, 49).Value = 22
$ws.Cells.Item(219, 60).Value = 60
$ws.Cells.Item(219, 63).Value = 90
$ws.Cells.Item(219, 65).Value = 66
$ws.Cells.Item(219, 66).Value = 18

# Row 220
$ws.Cells.Item(220, 2).Value = 7761
$ws.Cells.Item(220, 3).Value = 4953
$ws.Cells.Item(220, 4).Value = 4554
$ws.Cells.Item(220, 5).Value = 1825
$ws.Cells.Item(220, 6).Value = 2547
$ws.Cells.Item(220, 13).Value = 638
$ws.Cells.Item(220, 14).Value = 580
$ws.Cells.Item(220, 15).Value = 260
$ws.Cells.Item(220, 16).Value = 44
$ws.Cells.Item(220, 17).Value = 31
$ws.Cells.Item(220, 18).Value = 86
$ws.Cells.Item(220, 20).Value = 88
$ws.Cells.Item(220, 21).Value = 20
$ws.Cells.Item(220, 23).Value = 45
$ws.Cells.Item(220, 24).Value = 7
$ws.Cells.Item(220, 29).Value = 917
$ws.Cells.Item(220, 32).Value = 420
$ws.Cells.Item(220, 58).Value = 447
$ws.Cells.Item(220, 71).Value = 40

# Row 221
$ws.Cells.Item(221, 2).Value = 8160
$ws.Cells.Item(221, 3).Value = 5067
$ws.Cells.Item(221, 4).Value = 4598
$ws.Cells.Item(221, 5).Value = 1746
$ws.Cells.Item(221, 6).Value = 2616
$ws.Cells.Item(221, 10).Value = 31
$ws.Cells.Item(221, 22).Value = 2
$ws.Cells.Item(221, 24).Value = 28
$ws.Cells.Item(221, 28).Value = 2513
$ws.Cells.Item(221, 29).Value = 899
$ws.Cells.Item(221, 32).Value = 384
$ws.Cells.Item(221, 33).Value = 29
$ws.Cells.Item(221, 46).Value = 26
$ws.Cells.Item(221, 58).Value = 541
$ws.Cells.Item(221, 68).Value = 179
$ws.Cells.Item(221, 70).Value = 80
$ws.Cells.Item(221, 71).Value = 70
$ws.Cells.Item(221, 72).Value = 74

# Row 222
$ws.Cells.Item(222, 2).Value = 7621
$ws.Cells.Item(222, 3).Value = 5029
$ws.Cells.Item(222, 4).Value = 4663
$ws.Cells.Item(222, 5).Value = 1715
$ws.Cells.Item(222, 6).Value = 2667
$ws.Cells.Item(222, 28).Value = 2204
$ws.Cells.Item(222, 47).Value = 199
$ws.Cells.Item(222, 58).Value = 424
$ws.Cells.Item(222, 63).Value = 83
$ws.Cells.Item(222, 65).Value = 93
$ws.Cells.Item(222, 66).Value = 15
$ws.Cells.Item(222, 68).Value = 126
$ws.Cells.Item(222, 71).Value = 37

# Row 223
$ws.Cells.Item(223, 2).Value = 7453
$ws.Cells.Item(223, 3).Value = 4803
$ws.Cells.Item(223, 4).Value = 4375
$ws.Cells.Item(223, 5).Value = 1660
$ws.Cells.Item(223, 6).Value = 2585
$ws.Cells.Item(223, 28).Value = 2217
$ws.Cells.Item(223, 29).Value = 866
$ws.Cells.Item(223, 41).Value = 47
$ws.Cells.Item(223, 43).Value = 186
$ws.Cells.Item(223, 58).Value = 440
$ws.Cells.Item(223, 68).Value = 138
$ws.Cells.Item(223, 71).Value = 44
$ws.Cells.Item(223, 72).Value = 79

# --- Append new row 224 (series 01-07-2021) ---
# A224 holds a text label that looks like a date ("01-07-2021"); entering it
# directly would be auto-converted to a date serial by Excel's smart input,
# so build it as a text formula first, then paste back as a literal value
# (this keeps the cell a plain shared-string with no extra number format/style).
$aCell = $ws.Cells.Item(224, 1)
$aCell.Formula = '="01-07-2021"'
$aCell.Copy()
$aCell.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Cells.Item(224, 2).Value = 7943
$ws.Cells.Item(224, 3).Value = 5144
$ws.Cells.Item(224, 4).Value = 4584
$ws.Cells.Item(224, 5).Value = 1762
$ws.Cells.Item(224, 6).Value = 2586
$ws.Cells.Item(224, 7).Value = 347
$ws.Cells.Item(224, 8).Value = 35
$ws.Cells.Item(224, 9).Value = 61
$ws.Cells.Item(224, 10).Value = 28
$ws.Cells.Item(224, 11).Value = 84
$ws.Cells.Item(224, 12).Value = 6
$ws.Cells.Item(224, 13).Value = 419
$ws.Cells.Item(224, 14).Value = 346
$ws.Cells.Item(224, 15).Value = 2
$ws.Cells.Item(224, 16).Value = 100
$ws.Cells.Item(224, 17).Value = 10
$ws.Cells.Item(224, 18).Value = 0
$ws.Cells.Item(224, 19).Value = 40
$ws.Cells.Item(224, 20).Value = 0
$ws.Cells.Item(224, 21).Value = 0
$ws.Cells.Item(224, 22).Value = 0
$ws.Cells.Item(224, 23).Value = 56
$ws.Cells.Item(224, 24).Value = 6
$ws.Cells.Item(224, 25).Value = 25
$ws.Cells.Item(224, 26).Value = 4
$ws.Cells.Item(224, 27).Value = 13
$ws.Cells.Item(224, 28).Value = 2379
$ws.Cells.Item(224, 29).Value = 890
$ws.Cells.Item(224, 30).Value = 53
$ws.Cells.Item(224, 31).Value = 14
$ws.Cells.Item(224, 32).Value = 326
$ws.Cells.Item(224, 33).Value = 29
$ws.Cells.Item(224, 34).Value = 12
$ws.Cells.Item(224, 35).Value = 5
$ws.Cells.Item(224, 36).Value = 58
$ws.Cells.Item(224, 37).Value = 36
$ws.Cells.Item(224, 38).Value = 45
$ws.Cells.Item(224, 39).Value = 26
$ws.Cells.Item(224, 40).Value = 11
$ws.Cells.Item(224, 41).Value = 49
$ws.Cells.Item(224, 42).Value = 69
$ws.Cells.Item(224, 43).Value = 177
$ws.Cells.Item(224, 44).Value = 10
$ws.Cells.Item(224, 45).Value = 138
$ws.Cells.Item(224, 46).Value = 26
$ws.Cells.Item(224, 47).Value = 234
$ws.Cells.Item(224, 48).Value = 78
$ws.Cells.Item(224, 49).Value = 26
$ws.Cells.Item(224, 50).Value = 30
$ws.Cells.Item(224, 51).Value = 32
$ws.Cells.Item(224, 52).Value = 41
$ws.Cells.Item(224, 53).Value = 348
$ws.Cells.Item(224, 54).Value = 34
$ws.Cells.Item(224, 55).Value = 129
$ws.Cells.Item(224, 56).Value = 104
$ws.Cells.Item(224, 57).Value = 28
$ws.Cells.Item(224, 58).Value = 413
$ws.Cells.Item(224, 59).Value = 5
$ws.Cells.Item(224, 60).Value = 40
$ws.Cells.Item(224, 61).Value = 17
$ws.Cells.Item(224, 62).Value = 38
$ws.Cells.Item(224, 63).Value = 105
$ws.Cells.Item(224, 64).Value = 28
$ws.Cells.Item(224, 65).Value = 73
$ws.Cells.Item(224, 66).Value = 15
$ws.Cells.Item(224, 67).Value = 26
$ws.Cells.Item(224, 68).Value = 169
$ws.Cells.Item(224, 69).Value = 35
$ws.Cells.Item(224, 70).Value = 74
$ws.Cells.Item(224, 71).Value = 60
$ws.Cells.Item(224, 72).Value = 76
